# Auto commit at 2025-12-16 8:02:55.98
# Updates the Metrics sheet with refreshed totals and re-points the
# saved selections on the Metrics and today sheets. The today sheet's
# B/E/F columns are formulas that read straight from Metrics, so they
# (and A1's TODAY()-1) recalculate automatically once the source values
# change.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 204536.63
$metrics.Range("B3").Value = 175146.77000000002
$metrics.Range("B4").Value = 62756.109999999993
$metrics.Range("B5").Value = 8333
$metrics.Range("B6").Value = 5407243.7400000012
$metrics.Range("B7").Value = 4575499.7300000004
$metrics.Range("B8").Value = 1594712.9900000005
$metrics.Range("B9").Value = 211040
$metrics.Range("B10").Value = 33872624.729999989
$metrics.Range("B11").Value = 31850774.890000001
$metrics.Range("B12").Value = 11876435.029999996
$metrics.Range("B13").Value = 1308670

$metrics.Activate()
$metrics.Range("D12").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("H9").Select()
